$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 128, shifting existing rows 128-237 down to 129-238.
$ws.Rows.Item(128).Insert()

# Populate the newly inserted row 128 with the new data record.
$ws.Cells.Item(128, 1).Value = 11
$ws.Cells.Item(128, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(128, 3).Value = "Bíobío"
$ws.Cells.Item(128, 4).Value = 44944
$ws.Cells.Item(128, 5).Value = 8
$ws.Cells.Item(128, 6).Value = "Fruta"
$ws.Cells.Item(128, 7).Value = 100108
$ws.Cells.Item(128, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(128, 9).Value = 100108005
$ws.Cells.Item(128, 10).Value = "Piña"
$ws.Cells.Item(128, 11).Value = "Caramelo"
$ws.Cells.Item(128, 12).Value = "Segunda"
$ws.Cells.Item(128, 13).Value = 100
$ws.Cells.Item(128, 14).Value = 17000
$ws.Cells.Item(128, 15).Value = 18000
$ws.Cells.Item(128, 16).Value = 17500
$ws.Cells.Item(128, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(128, 18).Value = "Ecuador"
$ws.Cells.Item(128, 19).Value = 1250
$ws.Cells.Item(128, 20).Value = 14
